$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete last two data rows (old rows 8 and 9),
# bringing the table down to 7 rows (header + 6 data rows).
$ws.Rows.Item(9).Delete() | Out-Null
$ws.Rows.Item(8).Delete() | Out-Null

# Updated TPM-derived data for the remaining rows (2-7), columns A:T.
$data = @(
    @("FAPs","Efnb3","Ephb1","ECs",2,0.6666666666666666,0.5905326666666667,1.771598,0.504014252159901,0.5040142521599009,3,1,2.718682666666667,8.156048,0.5434637507613679,0.5434637507613679,1.605470924967111,14.449238324704,0.2739134759160057,0.2739134759160056),
    @("FAPs","Efnb3","Ephb1","MuSCs",2,0.6666666666666666,0.5905326666666667,1.771598,0.504014252159901,0.5040142521599009,3,1,2.283827,6.851481,0.4565362492386322,0.4565362492386321,1.348674448515333,12.138070036638,0.2301007762438954,0.2301007762438953),
    @("MuSCs","Efnb3","Ephb1","ECs",3,1,0.55299,1.65897,0.4719719281155831,0.4719719281155832,3,1,2.718682666666667,8.156048,0.5434637507613679,0.5434637507613679,1.50340432784,13.53063895056,0.2564996343077695,0.2564996343077696),
    @("MuSCs","Efnb3","Ephb1","MuSCs",3,1,0.55299,1.65897,0.4719719281155831,0.4719719281155832,3,1,2.283827,6.851481,0.4565362492386322,0.4565362492386321,1.26293349273,11.36640143457,0.2154722938078137,0.2154722938078137),
    @("Resolving-Mac","Efnb3","Ephb1","ECs",1,0.3333333333333333,0.028136,0.084408,0.0240138197245159,0.0240138197245159,3,1,2.718682666666667,8.156048,0.5434637507613679,0.5434637507613679,0.07649285550933332,0.688435699584,0.01305064053759273,0.01305064053759273),
    @("Resolving-Mac","Efnb3","Ephb1","MuSCs",1,0.3333333333333333,0.028136,0.084408,0.0240138197245159,0.0240138197245159,3,1,2.283827,6.851481,0.4565362492386322,0.4565362492386321,0.06425775647199999,0.578319808248,0.01096317918692317,0.01096317918692317)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value2 = $rowValues[$j]
    }
}
